$wb = $excel.ActiveWorkbook

# --- Step 1: duplicate the existing "effort" sheet; the copy becomes the new
#     "Effort R 1.0" sheet (placed before the original), while the original
#     remains and will be renamed to "Effort R 0.9" later. Using Copy() (instead
#     of Worksheets.Add()) means all existing formatting - column widths,
#     number formats, styles - are inherited exactly as-is.
$orig = $wb.Worksheets.Item("effort")
$orig.Copy($orig)

# The newly inserted copy is now worksheet #1 (placed right before $orig).
$newWs = $wb.Worksheets.Item(1)

# --- Step 2: make room for the new "Additional Effort [h]" column by
#     inserting a blank column at C. This pushes the former column C (task
#     descriptions) to column D, carrying its column width / bestFit setting
#     along with it.
$newWs.Columns.Item(3).Insert()

# --- Step 3: the new sheet only keeps the first 15 data rows (2-16); drop
#     everything below that.
$newWs.Range("A17:D58").EntireRow.Delete()

# --- Step 4: update the header row.
$newWs.Range("A1").Value = "Date"
$newWs.Range("B1").Value = "Effort [h]"
$newWs.Range("C1").Value = "Additional Effort [h]"
$newWs.Range("D1").Value = "Task"

# --- Step 5: write out the new data set (dates are Excel serial numbers, same
#     representation as the original workbook).
$data = @(
    @(41423, 2.5, 2.5,  "Revision of Makefile"),
    @(41424, 2,   $null, "Documentation Makefile changes. Concept for new sync objects"),
    @(41425, 2,   $null, "Concept of new sync objects"),
    @(41426, 0.75,$null, "Design of implementation new sync objects "),
    @(41430, 2,   $null, "Implementation of mutexes"),
    @(41431, 2,   $null, "Implementation of mutexes"),
    @(41432, 2,   2,     "Update Manual"),
    @(41432, 2.25,$null, "Implementation of mutexes"),
    @(41435, 2,   $null, "Implementation of mutexes: Basically done. No test case implemented yet, no testing done yet"),
    @(41436, 1.5, $null, "Implementation of tc11_mutex"),
    @(41439, 1.5, 2.5,  "Implementation of semaphores and first, very preliminary but successfuls tests"),
    @(41440, 2.5, $null, "Implementation tc12_queue"),
    @(41442, 2,   $null, "Implementation tc12_queue"),
    @(41443, 1,   3,     "Implementation tc12_queue"),
    @(41444, 1.5, 2.5,  "Design and implementation tc13_eventStates")
)

$r = 2
foreach ($row in $data) {
    $newWs.Range("A$r").Value = $row[0]
    $newWs.Range("B$r").Value = $row[1]
    if ($row[2] -ne $null) {
        $newWs.Range("C$r").Value = $row[2]
    }
    $newWs.Range("D$r").Value = $row[3]
    $r = $r + 1
}

# --- Step 6: rename the sheets. Re-fetch the original by its still-unique
#     name, since worksheet references captured before Copy()/Add() rebind by
#     position rather than identity.
$newWs.Name = "Effort R 1.0"
$oldWs = $wb.Worksheets.Item("effort")
$oldWs.Name = "Effort R 0.9"

# --- Step 7: tab/view bookkeeping - the old sheet's selection moves onto
#     column B, while the new sheet stays the selected/active tab (it must be
#     activated last so it remains the active one).
$oldWs.Range("B1:B1048576").Select()
$newWs.Range("D25").Select()
